# Commit message: "Fruta / hortaliza, semanal"
# The edit inserts one new weekly data row at row 208 of the sheet
# (shifting the former rows 208-298 down to 209-299), and grows the
# sheet dimension from A1:R298 to A1:R299 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 208; this automatically shifts all
# rows at/after 208 down by one and extends the used range.
$ws.Rows("208:208").Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A208").Value = 6
$ws.Range("B208").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C208").Value = "Metropolitana"
$ws.Range("D208").Value = 45009
$ws.Range("E208").Value = 13
$ws.Range("F208").Value = 100112001
$ws.Range("G208").Value = "Berenjena"
$ws.Range("H208").Value = "Sin especificar"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 300
$ws.Range("K208").Value = 5000
$ws.Range("L208").Value = 6000
$ws.Range("M208").Value = 5400
$ws.Range("N208").Value = "$/caja 60 unidades"
$ws.Range("O208").Value = "Región Metropolitana"
$ws.Range("P208").Value = 90
$ws.Range("Q208").Value = 60
$ws.Range("R208").Value = "Hortaliza"

# Keep the D column formatted the same as the rest of the date column
# (style already carried over from the row above by the Insert, but
# make sure it matches explicitly).
$ws.Range("D208").NumberFormat = $ws.Range("D209").NumberFormat
